$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set homework grade (5) for several students' assignments
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 5

$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 5

$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 5

$ws.Range("D17").Value = 5
$ws.Range("F17").Value = 5

$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 5

$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 5
$ws.Range("F22").Value = 5

$ws.Range("F24").Value = 5

$ws.Range("C30").Value = 5
$ws.Range("E30").Value = 5

# Update the active selection to match the final cursor position
$ws.Range("F18").Select()
